$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 1748
$ws.Range('L3').Value = 1771
$ws.Range('D4').Value = 1992
$ws.Range('E4').Value = 2046
$ws.Range('I4').Value = 1835
$ws.Range('K4').Value = 1757
$ws.Range('L4').Value = 501
$ws.Range('L5').Value = 107
$ws.Range('K6').Value = 9122
$ws.Range('L6').Value = 1652
$ws.Range('D7').Value = 28183
$ws.Range('E7').Value = 26051
$ws.Range('I7').Value = 26302
$ws.Range('K7').Value = 27546
$ws.Range('L7').Value = 5779

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L4').Value = 20
$ws.Range('L5').Value = 20
$ws.Range('L7').Value = 192
$ws.Range('L8').Value = 352
$ws.Range('L11').Value = 104
$ws.Range('L14').Value = 29
$ws.Range('L16').Value = 12
$ws.Range('L19').Value = 168
$ws.Range('L20').Value = 156
$ws.Range('L25').Value = 26
$ws.Range('L27').Value = 63
$ws.Range('L29').Value = 292
$ws.Range('L31').Value = 60
$ws.Range('L33').Value = 256
$ws.Range('L36').Value = 87
$ws.Range('L37').Value = 207
$ws.Range('L42').Value = 181
$ws.Range('L44').Value = 42
$ws.Range('L45').Value = 10
$ws.Range('L48').Value = 85
$ws.Range('L50').Value = 35
$ws.Range('L52').Value = 123
$ws.Range('K54').Value = 532
$ws.Range('L54').Value = 127
$ws.Range('L55').Value = 53
$ws.Range('L57').Value = 28
$ws.Range('D63').Value = 372
$ws.Range('E63').Value = 382
$ws.Range('I63').Value = 257
$ws.Range('K63').Value = 87
$ws.Range('L63').Value = 17
$ws.Range('L64').Value = 40
$ws.Range('L65').Value = 111
$ws.Range('L67').Value = 200
$ws.Range('L71').Value = 15
$ws.Range('L72').Value = 23
$ws.Range('L79').Value = 159
$ws.Range('L83').Value = 134
$ws.Range('L86').Value = 40
$ws.Range('L87').Value = 17
$ws.Range('L88').Value = 76
$ws.Range('L91').Value = 74
$ws.Range('L95').Value = 90
$ws.Range('L96').Value = 53
$ws.Range('L98').Value = 43
$ws.Range('L99').Value = 90
$ws.Range('D101').Value = 28183
$ws.Range('E101').Value = 26051
$ws.Range('I101').Value = 26302
$ws.Range('K101').Value = 27546
$ws.Range('L101').Value = 5779

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('L2').Value = 11
$ws.Range('L7').Value = 29

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L3').Value = 10
$ws.Range('L7').Value = 53

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 55
$ws.Range('L3').Value = 61
$ws.Range('L7').Value = 192

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L3').Value = 31
$ws.Range('L7').Value = 104

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 89
$ws.Range('L6').Value = 54

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L2').Value = 41
$ws.Range('L3').Value = 35
$ws.Range('L7').Value = 123

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 99
$ws.Range('L3').Value = 113
$ws.Range('L4').Value = 28
$ws.Range('L5').Value = 12
$ws.Range('L6').Value = 100
$ws.Range('L7').Value = 352

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L2').Value = 40
$ws.Range('L6').Value = 30
$ws.Range('L7').Value = 134

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 58
$ws.Range('L4').Value = 14
$ws.Range('L6').Value = 92
$ws.Range('L7').Value = 256

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L3').Value = 30
$ws.Range('L7').Value = 90

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 64
$ws.Range('L4').Value = 13
$ws.Range('L7').Value = 207

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 42
$ws.Range('L3').Value = 34
$ws.Range('L7').Value = 111

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L3').Value = 40
$ws.Range('L6').Value = 17
$ws.Range('L7').Value = 90

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L3').Value = 15
$ws.Range('L7').Value = 60

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 57
$ws.Range('L7').Value = 200

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K6').Value = 281
$ws.Range('L6').Value = 63
$ws.Range('K7').Value = 532
$ws.Range('L7').Value = 127

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L3').Value = 102
$ws.Range('L6').Value = 81
$ws.Range('L7').Value = 292

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L4').Value = 25
$ws.Range('L7').Value = 85

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 52
$ws.Range('L6').Value = 54
$ws.Range('L7').Value = 168

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('L3').Value = 11
$ws.Range('L7').Value = 42

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L2').Value = 46
$ws.Range('L3').Value = 46
$ws.Range('L7').Value = 181

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L2').Value = 21
$ws.Range('L7').Value = 53

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L2').Value = 27
$ws.Range('L3').Value = 24
$ws.Range('L7').Value = 74

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L3').Value = 51
$ws.Range('L4').Value = 12
$ws.Range('L6').Value = 34
$ws.Range('L7').Value = 159

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('L2').Value = 14
$ws.Range('L7').Value = 40

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 51
$ws.Range('L7').Value = 156

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L6').Value = 21
$ws.Range('L7').Value = 87

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('L3').Value = 13
$ws.Range('L7').Value = 26

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('L2').Value = 10
$ws.Range('L4').Value = 3
$ws.Range('L7').Value = 43

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('L4').Value = 3
$ws.Range('L7').Value = 35

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('L6').Value = 28
$ws.Range('L7').Value = 76

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('L2').Value = 5
$ws.Range('L7').Value = 20

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L3').Value = 24
$ws.Range('L7').Value = 63

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L4').Value = 23
$ws.Range('L7').Value = 40

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('L3').Value = 7
$ws.Range('L7').Value = 28

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('L3').Value = 9
$ws.Range('L7').Value = 15

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('L4').Value = 2
$ws.Range('L7').Value = 23

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range('L4').Value = 2
$ws.Range('L7').Value = 10

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('L2').Value = 5
$ws.Range('L7').Value = 20

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('L2').Value = 4
$ws.Range('L7').Value = 17

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('L6').Value = 9
$ws.Range('L7').Value = 12
